# Add the new "Saipan Beach" location as row 90 of the location-1 sheet.
# (Category, lat/long, Location, City, Country, YouTube video id)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row's values. The order below matches the order the
# corresponding shared strings were appended to xl/sharedStrings.xml in the
# original edit.
$ws.Range("B90").Value = "15.21426359540567, 145.71576906988687"
$ws.Range("F90").Value = "zFGugdfc8k4"
$ws.Range("D90").Value = "Saipan"
$ws.Range("C90").Value = "Saipan Beach"
$ws.Range("A90").Value = "LIVE, SEA, BEACH"
$ws.Range("E90").Value = "USA"

# Match the formatting used by the rest of the table: columns A/C/D/E carry
# the thin left/right border style used throughout the sheet, while B and F
# stay unformatted (same as every other row).
$ws.Range("A89").Copy()
$ws.Range("A90").PasteSpecial(-4122)

$ws.Range("C89").Copy()
$ws.Range("C90:D90").PasteSpecial(-4122)

$ws.Range("E89").Copy()
$ws.Range("E90").PasteSpecial(-4122)

# Reflect where the user ended up after typing the new row.
$ws.Range("A91").Select() | Out-Null
